$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: status text for zh-cn / de-de now reflects that the
# handback generation completed ("Ready for handoff" -> "Handed back:
# in sync with en-US"), and the two locale-status columns are widened
# to fit the longer text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: the handback report now fills in the Latest Target File
# (hyperlinked, same as the source-file link), Latest Handback File and
# Latest Handback DateTime columns, and widens the columns that hold
# longer file names.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsZh.Range("I2").Value = "3d71b4c5-947a-429a-9c10-4e23a09dd6d6.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96ae5ea01a86954fb4fcbd661e640428f397117e/e2e/3d71b4c5-947a-429a-9c10-4e23a09dd6d6.md", "", "3d71b4c5-947a-429a-9c10-4e23a09dd6d6.md", "3d71b4c5-947a-429a-9c10-4e23a09dd6d6.md")
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276

$wsZh.Range("J2").Value = "3d71b4c5-947a-429a-9c10-4e23a09dd6d6.dc1311b846f9dd62cbf907a065a0b0c12964926f.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-21 07:04:40"

# ---------------------------------------------------------------------
# de-de sheet: same shape of update as zh-cn, but with the de-de xlf
# file name and its own handback timestamp.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Range("I2").Value = "3d71b4c5-947a-429a-9c10-4e23a09dd6d6.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96ae5ea01a86954fb4fcbd661e640428f397117e/e2e/3d71b4c5-947a-429a-9c10-4e23a09dd6d6.md", "", "3d71b4c5-947a-429a-9c10-4e23a09dd6d6.md", "3d71b4c5-947a-429a-9c10-4e23a09dd6d6.md")
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276

$wsDe.Range("J2").Value = "3d71b4c5-947a-429a-9c10-4e23a09dd6d6.dc1311b846f9dd62cbf907a065a0b0c12964926f.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-21 07:04:46"
